$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2641.1667
$ws.Range("I32").Value = 2000
$ws.Range("J32").Value = 2769.4
$ws.Range("K32").Value = 2000
$ws.Range("L32").Value = 2769.4
$ws.Range("M32").Value = -1674
$ws.Range("N32").Value = -3421.4
$ws.Range("H43").Value = 11682
$ws.Range("I43").Value = 8999
$ws.Range("K43").Value = 8999
$ws.Range("M43").Value = -8930
$ws.Range("H80").Value = 1030.0385
$ws.Range("I80").Value = 1132.4
$ws.Range("K80").Value = 3397.2
$ws.Range("M80").Value = -2399.2
$ws.Range("H83").Value = 1030.0385
$ws.Range("I83").Value = 1132.4
$ws.Range("K83").Value = 10191.6
$ws.Range("M83").Value = -5199.6
$ws.Range("H128").Value = 74999.5
$ws.Range("J128").Value = 74999.5
$ws.Range("L128").Value = 74999.5
$ws.Range("N128").Value = -84959.5
$ws.Range("H131").Value = 9699.143
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 9699.143
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 29097.429
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -39177.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 11000
$ws.Range("J2").Value = 11000
$ws.Range("L2").Value = 11000
$ws.Range("N2").Value = -11226
$ws.Range("H45").Value = 1946.6923
$ws.Range("I45").Value = 1734.7273
$ws.Range("K45").Value = 1734.7273
$ws.Range("M45").Value = -1357.7273
$ws.Range("H61").Value = 5092.353
$ws.Range("I61").Value = 5092.353
$ws.Range("K61").Value = 5092.353
$ws.Range("M61").Value = -4880.353
$ws.Range("H74").Value = 1621.6
$ws.Range("I74").Value = 1601.7368
$ws.Range("K74").Value = 1601.7368
$ws.Range("M74").Value = -727.7367999999999
$ws.Range("H77").Value = 1621.6
$ws.Range("I77").Value = 1601.7368
$ws.Range("K77").Value = 8008.683999999999
$ws.Range("M77").Value = -3640.683999999999
$ws.Range("H97").Value = 471.38095
$ws.Range("J97").Value = 597.8333
$ws.Range("L97").Value = 597.8333
$ws.Range("N97").Value = -1589.8333
$ws.Range("H116").Value = 11000
$ws.Range("J116").Value = 11000
$ws.Range("L116").Value = 11000
$ws.Range("N116").Value = -15588
$ws.Range("H122").Value = 2724.375
$ws.Range("I122").Value = 2492
$ws.Range("J122").Value = 2956.75
$ws.Range("K122").Value = 7476
$ws.Range("L122").Value = 8870.25
$ws.Range("M122").Value = -5026
$ws.Range("N122").Value = -13770.25
$ws.Range("H132").Value = 1462.2812
$ws.Range("I132").Value = 1462.2812
$ws.Range("K132").Value = 4386.8436
$ws.Range("M132").Value = -1856.8436
$ws.Range("H136").Value = 5092.353
$ws.Range("I136").Value = 5092.353
$ws.Range("K136").Value = 15277.059
$ws.Range("M136").Value = -12727.059

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 11000
$ws.Range("J3").Value = 11000
$ws.Range("L3").Value = 11000
$ws.Range("N3").Value = -11228
$ws.Range("H82").Value = 15982.75
$ws.Range("I82").Value = 15982.75
$ws.Range("K82").Value = 15982.75
$ws.Range("M82").Value = -15599.75
$ws.Range("H85").Value = 15982.75
$ws.Range("I85").Value = 15982.75
$ws.Range("K85").Value = 15982.75
$ws.Range("M85").Value = -14656.75
$ws.Range("H94").Value = 2165.4546
$ws.Range("I94").Value = 1424.4445
$ws.Range("J94").Value = 5500
$ws.Range("K94").Value = 1424.4445
$ws.Range("L94").Value = 5500
$ws.Range("M94").Value = -973.4445000000001
$ws.Range("N94").Value = -6402
$ws.Range("H134").Value = 2738.8215
$ws.Range("I134").Value = 2662.2632
$ws.Range("K134").Value = 7986.7896
$ws.Range("M134").Value = -5451.7896

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2325.75
$ws.Range("I58").Value = 2406
$ws.Range("K58").Value = 2406
$ws.Range("M58").Value = -2203
$ws.Range("H86").Value = 4754.857
$ws.Range("I86").Value = 5395.25
$ws.Range("J86").Value = 3901
$ws.Range("K86").Value = 5395.25
$ws.Range("L86").Value = 3901
$ws.Range("M86").Value = -4272.25
$ws.Range("N86").Value = -6147
$ws.Range("H89").Value = 4754.857
$ws.Range("I89").Value = 5395.25
$ws.Range("J89").Value = 3901
$ws.Range("K89").Value = 26976.25
$ws.Range("L89").Value = 19505
$ws.Range("M89").Value = -21360.25
$ws.Range("N89").Value = -30737
$ws.Range("H132").Value = 2056.9
$ws.Range("I132").Value = 2056.9
$ws.Range("K132").Value = 6170.700000000001
$ws.Range("M132").Value = -3640.700000000001
$ws.Range("H136").Value = 2325.75
$ws.Range("I136").Value = 2406
$ws.Range("K136").Value = 7218
$ws.Range("M136").Value = -4668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 8069728
$ws.Range("I4").Value = 1697990.6
$ws.Range("J4").Value = 62866668
$ws.Range("K4").Value = 5093971.800000001
$ws.Range("L4").Value = 188600004
$ws.Range("M4").Value = -5093859.800000001
$ws.Range("N4").Value = -188600228
$ws.Range("H44").Value = 87.5
$ws.Range("I44").Value = 87.5
$ws.Range("K44").Value = 262.5
$ws.Range("M44").Value = 135.5
$ws.Range("H52").Value = 1494.4
$ws.Range("J52").Value = 1494.4
$ws.Range("L52").Value = 4483.200000000001
$ws.Range("N52").Value = -5015.200000000001
$ws.Range("H103").Value = 713.55554
$ws.Range("J103").Value = 1742.6666
$ws.Range("L103").Value = 5227.9998
$ws.Range("N103").Value = -6985.9998
$ws.Range("H114").Value = 1088.1111
$ws.Range("J114").Value = 1461.75
$ws.Range("L114").Value = 4385.25
$ws.Range("N114").Value = -10893.25
$ws.Range("H129").Value = 3554.6
$ws.Range("J129").Value = 4716.1816
$ws.Range("L129").Value = 14148.5448
$ws.Range("N129").Value = -24148.5448

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 3500
$ws.Range("J27").Value = 5000
$ws.Range("L27").Value = 5000
$ws.Range("N27").Value = -5332
$ws.Range("H46").Value = 26856.572
$ws.Range("I46").Value = 9333
$ws.Range("J46").Value = 39999.25
$ws.Range("K46").Value = 9333
$ws.Range("L46").Value = 39999.25
$ws.Range("M46").Value = -9177
$ws.Range("N46").Value = -40311.25
$ws.Range("H93").Value = 26049.4
$ws.Range("J93").Value = 26049.4
$ws.Range("L93").Value = 26049.4
$ws.Range("N93").Value = -29793.4
$ws.Range("H122").Value = 3388.8667
$ws.Range("I122").Value = 2752.25
$ws.Range("J122").Value = 3620.3635
$ws.Range("K122").Value = 8256.75
$ws.Range("L122").Value = 10861.0905
$ws.Range("M122").Value = -5806.75
$ws.Range("N122").Value = -15761.0905
$ws.Range("H132").Value = 2454.5557
$ws.Range("I132").Value = 2454.5557
$ws.Range("K132").Value = 7363.6671
$ws.Range("M132").Value = -4833.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3631.3333
$ws.Range("I7").Value = 1811.7142
$ws.Range("K7").Value = 1811.7142
$ws.Range("M7").Value = -1699.7142
$ws.Range("H40").Value = 3719.2334
$ws.Range("I40").Value = 2678.9
$ws.Range("K40").Value = 2678.9
$ws.Range("M40").Value = -2542.9
$ws.Range("H122").Value = 9589
$ws.Range("I122").Value = 9611.25
$ws.Range("K122").Value = 28833.75
$ws.Range("M122").Value = -26383.75
$ws.Range("H126").Value = 3631.3333
$ws.Range("I126").Value = 1811.7142
$ws.Range("K126").Value = 5435.142599999999
$ws.Range("M126").Value = -2965.142599999999
$ws.Range("H133").Value = 49000
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3429.875
$ws.Range("I132").Value = 3164.6538
$ws.Range("K132").Value = 9493.9614
$ws.Range("M132").Value = -6963.9614
